$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "https://www.nordea.se/"
$ws.Range("A7").Value = "https://www.aftonbladet.se/"

$ws.Range("A7").Select()
